$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.290778040885925
$ws.Range("B1").Value = 1.969690203666687
$ws.Range("C1").Value = 2.691423654556274
$ws.Range("D1").Value = 3.727025985717773
$ws.Range("E1").Value = 1.046867847442627
